# PlayerPerformance_5952.xlsx update
# 1. Insert a new "Player Info" sheet at the front of the workbook with
#    player identity details.
# 2. On the "ODI Batting" sheet, rename MATCH_CARD_LINK -> MATCH_CODE and
#    replace the scorecard URLs with the bare match codes.
# 3. On the "ODI Bowling" sheet, do the same rename/replacement.

$wb = $excel.ActiveWorkbook

$battingSheet = $wb.Worksheets.Item("ODI Batting")

# --- 1. Add the new "Player Info" sheet, placed before "ODI Batting" ---
$playerInfo = $wb.Worksheets.Add($battingSheet)
$playerInfo.Name = "Player Info"

# NOTE: inserting a sheet shifts worksheet references that were obtained
# by index/position, so re-fetch every sheet we still need by name before
# touching it again.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")
$playerInfo = $wb.Worksheets.Item("Player Info")

# Copy the header formatting (bold, bordered, centered) from an existing
# header cell so the new sheet's header row matches the workbook style.
$battingSheet.Range("A1").Copy()
$playerInfo.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

# Values are kept as plain text to mirror the original inline-string data.
$playerInfo.Range("A2:D2").NumberFormat = "@"
$playerInfo.Range("A2").Value = "5952"
$playerInfo.Range("B2").Value = "Kevin Sinclair"
$playerInfo.Range("C2").Value = "Right Handed"
$playerInfo.Range("D2").Value = "Right Arm Off Break"

$playerInfo.Range("A1").Select()

# --- 2. "ODI Batting": MATCH_CARD_LINK -> MATCH_CODE ---
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingSheet.Range("D2:D4").NumberFormat = "@"
$battingSheet.Range("D2").Value = "4636"
$battingSheet.Range("D3").Value = "4639"
$battingSheet.Range("D4").Value = "4642"

# --- 3. "ODI Bowling": MATCH_CARD_LINK -> MATCH_CODE ---
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingSheet.Range("B2:B4").NumberFormat = "@"
$bowlingSheet.Range("B2").Value = "4636"
$bowlingSheet.Range("B3").Value = "4639"
$bowlingSheet.Range("B4").Value = "4642"
